$d = $word.ActiveDocument

# The document originally had two extra class definitions ("Clase genero"
# and "Clase Director") that had been appended as their own paragraphs
# after the "Clase Pelicula" paragraph. This reverts that addition: it
# deletes both paragraphs' text (and the paragraph mark that separated
# them from the preceding paragraph), while preserving the trailing
# manual line break that belonged to the last ("Clase Director")
# paragraph. That break ends up appended right after the line break
# already at the end of the "Clase Pelicula" paragraph - exactly as the
# document looked before those two paragraphs were added.

$vt = [char]11
$cr = [char]13

$findText = $cr + "Clase genero" + $vt + "id(int), nombre(string)" + $cr + "Clase Director" + $vt + "id(int), nombre(string), apellido(string)"

$found = $d.Content.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

if (-not $found) {
    throw "Could not find the 'Clase genero' / 'Clase Director' paragraphs to remove."
}
